$wb = $excel.ActiveWorkbook

# Map of worksheet name -> list of (cell, newValue) updates to the "想去人数" (F) column
$updates = @{
    "展览" = @(
        @{ Cell = "F2";  Value = 851 },
        @{ Cell = "F3";  Value = 983 },
        @{ Cell = "F4";  Value = 773 },
        @{ Cell = "F5";  Value = 864 },
        @{ Cell = "F6";  Value = 435 },
        @{ Cell = "F8";  Value = 150 },
        @{ Cell = "F9";  Value = 1267 },
        @{ Cell = "F10"; Value = 696 },
        @{ Cell = "F11"; Value = 408 },
        @{ Cell = "F15"; Value = 882 },
        @{ Cell = "F18"; Value = 371 },
        @{ Cell = "F20"; Value = 576 },
        @{ Cell = "F21"; Value = 136 },
        @{ Cell = "F23"; Value = 33 },
        @{ Cell = "F24"; Value = 936 }
    )
    "演出" = @(
        @{ Cell = "F2";  Value = 338 },
        @{ Cell = "F5";  Value = 641 },
        @{ Cell = "F10"; Value = 27 },
        @{ Cell = "F11"; Value = 109 }
    )
    "本地生活" = @(
        @{ Cell = "F2";  Value = 374 }
    )
    "全部类型" = @(
        @{ Cell = "F2";  Value = 374 },
        @{ Cell = "F3";  Value = 338 },
        @{ Cell = "F4";  Value = 851 },
        @{ Cell = "F5";  Value = 983 },
        @{ Cell = "F6";  Value = 773 },
        @{ Cell = "F7";  Value = 864 },
        @{ Cell = "F8";  Value = 435 },
        @{ Cell = "F10"; Value = 150 },
        @{ Cell = "F11"; Value = 1267 },
        @{ Cell = "F12"; Value = 696 },
        @{ Cell = "F15"; Value = 408 },
        @{ Cell = "F17"; Value = 641 },
        @{ Cell = "F20"; Value = 882 },
        @{ Cell = "F24"; Value = 371 },
        @{ Cell = "F28"; Value = 576 },
        @{ Cell = "F30"; Value = 27 },
        @{ Cell = "F31"; Value = 109 },
        @{ Cell = "F32"; Value = 109 },
        @{ Cell = "F33"; Value = 136 },
        @{ Cell = "F35"; Value = 33 },
        @{ Cell = "F36"; Value = 936 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Range($u.Cell).Value = $u.Value
    }
}
